$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 1000
$ws.Range("J31").Value = 1000
$ws.Range("L31").Value = 3000
$ws.Range("N31").Value = -3460

$ws.Range("H64").Value = 3911.68
$ws.Range("I64").Value = 3673.75
$ws.Range("J64").Value = 4334.6665
$ws.Range("K64").Value = 3673.75
$ws.Range("L64").Value = 4334.6665
$ws.Range("M64").Value = -3425.75
$ws.Range("N64").Value = -4830.6665

$ws.Range("H67").Value = 3911.68
$ws.Range("I67").Value = 3673.75
$ws.Range("J67").Value = 4334.6665
$ws.Range("K67").Value = 3673.75
$ws.Range("L67").Value = 4334.6665
$ws.Range("M67").Value = -2815.75
$ws.Range("N67").Value = -6050.6665

$ws.Range("H69").Value = 4517.115
$ws.Range("I69").Value = 4878.5835
$ws.Range("J69").Value = 4207.2856
$ws.Range("K69").Value = 14635.7505
$ws.Range("L69").Value = 12621.8568
$ws.Range("M69").Value = -13761.7505
$ws.Range("N69").Value = -14369.8568

$ws.Range("H72").Value = 4517.115
$ws.Range("I72").Value = 4878.5835
$ws.Range("J72").Value = 4207.2856
$ws.Range("K72").Value = 43907.2515
$ws.Range("L72").Value = 37865.5704
$ws.Range("M72").Value = -39539.2515
$ws.Range("N72").Value = -46601.5704

$ws.Range("H76").Value = 3485.147
$ws.Range("I76").Value = 3136.4583
$ws.Range("K76").Value = 3136.4583
$ws.Range("M76").Value = -2821.4583

$ws.Range("H79").Value = 3485.147
$ws.Range("I79").Value = 3136.4583
$ws.Range("K79").Value = 3136.4583
$ws.Range("M79").Value = -2044.4583

$ws.Range("H141").Value = 2716.75
$ws.Range("I141").Value = 2321.318
$ws.Range("J141").Value = 4166.6665
$ws.Range("K141").Value = 6963.954000000001
$ws.Range("L141").Value = 12499.9995
$ws.Range("M141").Value = -1783.954000000001
$ws.Range("N141").Value = -22859.9995


# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4963.8643
$ws.Range("I32").Value = 4278.579
$ws.Range("J32").Value = 24494.5
$ws.Range("K32").Value = 4278.579
$ws.Range("L32").Value = 24494.5
$ws.Range("M32").Value = -3991.579
$ws.Range("N32").Value = -25068.5

$ws.Range("H61").Value = 4658.3823
$ws.Range("I61").Value = 3654.92
$ws.Range("K61").Value = 3654.92
$ws.Range("M61").Value = -3442.92

$ws.Range("H63").Value = 4499.6665
$ws.Range("I63").Value = 4499.6665
$ws.Range("K63").Value = 4499.6665
$ws.Range("M63").Value = -3813.6665

$ws.Range("H66").Value = 4499.6665
$ws.Range("I66").Value = 4499.6665
$ws.Range("K66").Value = 22498.3325
$ws.Range("M66").Value = -19066.3325

$ws.Range("H74").Value = 1600.35
$ws.Range("I74").Value = 1600.35
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 1600.35
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -726.3499999999999
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 1600.35
$ws.Range("I77").Value = 1600.35
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 8001.75
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -3633.75
$ws.Range("N77").ClearContents()

$ws.Range("H102").Value = 2283.1667
$ws.Range("I102").Value = 2037.7778
$ws.Range("J102").Value = 2528.5557
$ws.Range("K102").Value = 2037.7778
$ws.Range("L102").Value = 2528.5557
$ws.Range("M102").Value = -415.7778000000001
$ws.Range("N102").Value = -5772.5557

$ws.Range("H132").Value = 9172.105
$ws.Range("I132").Value = 3362.2856
$ws.Range("J132").Value = 12561.167
$ws.Range("K132").Value = 10086.8568
$ws.Range("L132").Value = 37683.501
$ws.Range("M132").Value = -7556.856800000001
$ws.Range("N132").Value = -42743.501

$ws.Range("H133").Value = 20521.867
$ws.Range("J133").Value = 20521.867
$ws.Range("L133").Value = 20521.867
$ws.Range("N133").Value = -25581.867

$ws.Range("H136").Value = 4658.3823
$ws.Range("I136").Value = 3654.92
$ws.Range("K136").Value = 10964.76
$ws.Range("M136").Value = -8414.76


# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 29000
$ws.Range("J62").Value = 28000
$ws.Range("L62").Value = 28000
$ws.Range("N62").Value = -29372

$ws.Range("H65").Value = 29000
$ws.Range("J65").Value = 28000
$ws.Range("L65").Value = 84000
$ws.Range("N65").Value = -90864

$ws.Range("H86").Value = 2275.5715
$ws.Range("I86").Value = 2620.2
$ws.Range("J86").Value = 1414
$ws.Range("K86").Value = 2620.2
$ws.Range("L86").Value = 1414
$ws.Range("M86").Value = -1497.2
$ws.Range("N86").Value = -3660

$ws.Range("H89").Value = 2275.5715
$ws.Range("I89").Value = 2620.2
$ws.Range("J89").Value = 1414
$ws.Range("K89").Value = 13101
$ws.Range("L89").Value = 7070
$ws.Range("M89").Value = -7485
$ws.Range("N89").Value = -18302

$ws.Range("H105").Value = 3558.359
$ws.Range("I105").Value = 2942.6538
$ws.Range("J105").Value = 4789.769
$ws.Range("K105").Value = 2942.6538
$ws.Range("L105").Value = 4789.769
$ws.Range("M105").Value = -1195.6538
$ws.Range("N105").Value = -8283.769

$ws.Range("H134").Value = 4881.8203
$ws.Range("I134").Value = 5996.7144
$ws.Range("J134").Value = 2043.909
$ws.Range("K134").Value = 17990.1432
$ws.Range("L134").Value = 6131.727000000001
$ws.Range("M134").Value = -15455.1432
$ws.Range("N134").Value = -11201.727


# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3379.72
$ws.Range("I31").Value = 2523.75
$ws.Range("J31").Value = 4901.4443
$ws.Range("K31").Value = 2523.75
$ws.Range("L31").Value = 4901.4443
$ws.Range("M31").Value = -2228.75
$ws.Range("N31").Value = -5491.4443

$ws.Range("H34").Value = 3379.72
$ws.Range("I34").Value = 2523.75
$ws.Range("J34").Value = 4901.4443
$ws.Range("K34").Value = 2523.75
$ws.Range("L34").Value = 4901.4443
$ws.Range("M34").Value = -2321.75
$ws.Range("N34").Value = -5305.4443

$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()

$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

$ws.Range("H58").Value = 3250136.5
$ws.Range("I58").Value = 10103476
$ws.Range("J58").Value = 3817.5264
$ws.Range("K58").Value = 10103476
$ws.Range("L58").Value = 3817.5264
$ws.Range("M58").Value = -10103273
$ws.Range("N58").Value = -4223.526400000001

$ws.Range("H136").Value = 3250136.5
$ws.Range("I136").Value = 10103476
$ws.Range("J136").Value = 3817.5264
$ws.Range("K136").Value = 30310428
$ws.Range("L136").Value = 11452.5792
$ws.Range("M136").Value = -30307878
$ws.Range("N136").Value = -16552.5792


# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 5000
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 5000
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 15000
$ws.Range("M25").ClearContents()
$ws.Range("N25").Value = -15338

$ws.Range("H30").Value = 5000
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 5000
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 15000
$ws.Range("M30").ClearContents()
$ws.Range("N30").Value = -15204

$ws.Range("H44").Value = 1025
$ws.Range("I44").Value = 966.6667
$ws.Range("K44").Value = 2900.0001
$ws.Range("M44").Value = -2502.0001

$ws.Range("H123").Value = 3737.7144
$ws.Range("J123").Value = 3920.6316
$ws.Range("L123").Value = 11761.8948
$ws.Range("N123").Value = -16661.8948


# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()

$ws.Range("H80").Value = 6386.6665
$ws.Range("I80").Value = 12320
$ws.Range("J80").Value = 3420
$ws.Range("K80").Value = 12320
$ws.Range("L80").Value = 3420
$ws.Range("M80").Value = -11322
$ws.Range("N80").Value = -5416

$ws.Range("H83").Value = 6386.6665
$ws.Range("I83").Value = 12320
$ws.Range("J83").Value = 3420
$ws.Range("K83").Value = 61600
$ws.Range("L83").Value = 17100
$ws.Range("M83").Value = -56608
$ws.Range("N83").Value = -27084

$ws.Range("H132").Value = 3749.5
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 3749.5
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 11248.5
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -16308.5

$ws.Range("H134").Value = 37877.8
$ws.Range("J134").Value = 37877.8
$ws.Range("L134").Value = 113633.4
$ws.Range("N134").Value = -118703.4


# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4062.5
$ws.Range("I7").Value = 3612.5
$ws.Range("J7").Value = 4512.5
$ws.Range("K7").Value = 3612.5
$ws.Range("L7").Value = 4512.5
$ws.Range("M7").Value = -3500.5
$ws.Range("N7").Value = -4736.5

$ws.Range("H63").Value = 44542.5
$ws.Range("J63").Value = 44542.5
$ws.Range("L63").Value = 44542.5
$ws.Range("N63").Value = -46040.5

$ws.Range("H66").Value = 44542.5
$ws.Range("J66").Value = 44542.5
$ws.Range("L66").Value = 133627.5
$ws.Range("N66").Value = -141115.5

$ws.Range("H126").Value = 4062.5
$ws.Range("I126").Value = 3612.5
$ws.Range("J126").Value = 4512.5
$ws.Range("K126").Value = 10837.5
$ws.Range("L126").Value = 13537.5
$ws.Range("M126").Value = -8367.5
$ws.Range("N126").Value = -18477.5

$ws.Range("H132").Value = 4564.25
$ws.Range("I132").Value = 3880.6667
$ws.Range("K132").Value = 11642.0001
$ws.Range("M132").Value = -9112.000100000001


# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()

$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()

$ws.Range("H132").Value = 2574.682
$ws.Range("I132").Value = 1875.6666
$ws.Range("K132").Value = 5626.9998
$ws.Range("M132").Value = -3096.9998

